# SCR 12897: modify and add DIM_Source values
#  - rename "CMS Customer Call Listening" -> "CMS Reported Item"
#  - add a new "Internal CCO Reporting" sub-source for both Direct (118) and Indirect (218)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DIM_Source")

# Rename "CMS Customer Call Listening" -> "CMS Reported Item" wherever it appears
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    if ($ws.Cells.Item($r, 3).Text -eq "CMS Customer Call Listening") {
        $ws.Cells.Item($r, 3).Value = "CMS Reported Item"
    }
}

# Insert a new "Internal CCO Reporting" row for the Direct group (SourceID 118),
# right after the last existing Direct row (117 / BCC Security and Privacy Incident Coaching)
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = 118
$ws.Cells.Item(21, 2).Value = "Direct"
$ws.Cells.Item(21, 3).Value = "Internal CCO Reporting"

# Append a new "Internal CCO Reporting" row for the Indirect group (SourceID 218)
# as the new last row of the table
$newLastRow = $ws.UsedRange.Rows.Count + 1
$ws.Cells.Item($newLastRow, 1).Value = 218
$ws.Cells.Item($newLastRow, 2).Value = "Indirect"
$ws.Cells.Item($newLastRow, 3).Value = "Internal CCO Reporting"

# Make DIM_Source the active sheet/tab, with C39 selected, matching the saved view state
$ws.Activate()
$ws.Range("C" + $newLastRow).Select()
